$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous entry block (rows 15-19) onto the new
# entry block (rows 21-25) row by row, so the new block picks up the same
# wrap-text / left-top alignment style without introducing stray styles.
$ws.Range("H15:U15").Copy()
$ws.Range("H21:U21").PasteSpecial(-4122)
$ws.Range("H16:U16").Copy()
$ws.Range("H22:U22").PasteSpecial(-4122)
$ws.Range("H17:U17").Copy()
$ws.Range("H23:U23").PasteSpecial(-4122)
$ws.Range("H18:U18").Copy()
$ws.Range("H24:U24").PasteSpecial(-4122)
$ws.Range("H19:U19").Copy()
$ws.Range("H25:U25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New timesheet entry: 2017-03-31, 12h30 -> 15h00, 150 minutes
$ws.Range("B21").Value = 42825
$ws.Range("B21").NumberFormat = $ws.Range("B15").NumberFormat

$ws.Range("D21").Value = "12h30"
$ws.Range("E21").Value = "15h00"
$ws.Range("F21").Value = 150

$ws.Range("H21").Value = "Implémentation d'un logger pour avoir des informations et feedback plus parlant. Ajout de la copie de fichier selon le changement ainsi que la suppression de fichier lorsqu'il disparait. Ajout de la gestion d'Erreur de fichier. Ajout de la répétition des opérations lorsd'une erreur jusqu'a un nombre de répétition. "

# Merge the description block like the other entries
$ws.Range("H21:U25").Merge() | Out-Null

# Select the newly merged range, matching the recorded selection state
$ws.Range("H21:U25").Select() | Out-Null
